$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Daily practice")
$ws3 = $wb.Worksheets.Item("Learnings")

# Learnings sheet: new note about GCD property (string index 113 allocated first)
$ws3.Range("A5").Value = "GCD if done on array, is always decreasing.`ni.e, arr = {a1, a2, a3}`ngcd(a1, a2) -> g1`ngcd(a1,a2,a3) -> gcd(g1,a3) -> g2`nThen g2 <= g1`nSince, gcd(a,b) <= min(a,b)`n"
$ws3.Range("A5").WrapText = $true
$ws3.Rows.Item(5).RowHeight = 100.8

# Daily practice sheet: row 27 (Serval and Mocha's array, CF1789A)
$ws1.Range("A26:E26").Copy()
$ws1.Range("A27:E27").PasteSpecial(-4122)
$ws1.Range("A27").Value = "29th May, 2025 (HAD AN INTERVIEW)"
$ws1.Range("C27").Value = "https://codeforces.com/problemset/problem/1789/A"
$ws1.Range("D27").Value = "Difficult (since in-depth math of gcd)"
$ws1.Range("B27").Value = "Serval and Mocha's array"
$ws1.Range("E27").Value = "I initially considered sorting the array and then calculating gcd of every prefix . But that fails for :`n5`n1261 227821 143 4171 1941`nSo, I had to look at the logic and the interesting proerpty I learnt is:`nGCD if done on array, is always decreasing.`ni.e, arr = {a1, a2, a3}`ngcd(a1, a2) -> g1`ngcd(a1,a2,a3) -> gcd(g1,a3) -> g2`nThen g2 <= g1`nSince, gcd(a,b) <= min(a,b)`nThus, the idea is if we find any two elements with gcd <= 2, then we get the solution.`nSince the first prefix is ought to be <= 2 (since the array is of length 2)`nGCD is always decreasing, so if we find such two elements, the next prefixes will have gcd <=2 so it will always be a GOOD array"
$ws1.Range("A27:E27").RowHeight = 230.4

# row 28 (One and two, CF1788A)
$ws1.Range("A27:E27").Copy()
$ws1.Range("A28:E28").PasteSpecial(-4122)
$ws1.Range("A28").Value = "30th May, 2025"
$ws1.Range("B28").Value = "One and two"
$ws1.Range("C28").Value = "https://codeforces.com/problemset/problem/1788/A"
$ws1.Range("D28").Value = "Medium"
$ws1.Range("E28").Value = "All lied in the number of 2s in the array.`nIf all are 1, then k=1 is always correct, since the product till k will be equal to the product of elements from k+1 to n and all will be 1`nIf there are 2 present, we count the number of 2s.`nIf even number of 2s present, then the num of 2s/2 position 2 is the value of k`nElse, we cannot get k hence return -1"
$ws1.Range("A28:E28").RowHeight = 100.8

# row 29 (Make it beautiful, CF1783A)
$ws1.Range("A28:E28").Copy()
$ws1.Range("A29:E29").PasteSpecial(-4122)
$ws1.Range("A29").Value = "30th May, 2025"
$ws1.Range("B29").Value = "Make it beautiful"
$ws1.Range("C29").Value = "https://codeforces.com/problemset/problem/1783/A"
$ws1.Range("D29").Value = "Easy"
$ws1.Range("E29").Value = "I thought of a solution where I essentially make the array descending and am assuming that will always make the array beautiful because we cannot get any element that will be equal to the sum of elements before it since all elements before it will be larger than it !`nBut this did not work `nI found the case myself with some analysis.`nSince i was checkign for the first 2 elemetns and if they are same just exchanging the 2nd and 3rd elemetns. This process was incrorrect since it would fail for :`n3 3 3 2`nSo, I am now finding the next non similar element`nAnd that works !"
$ws1.Range("A29:E29").RowHeight = 201.6

# row 30 (Everybody likes good arrays, CF1777A)
$ws1.Range("A29:E29").Copy()
$ws1.Range("A30:E30").PasteSpecial(-4122)
$ws1.Range("A30").Value = "30th May, 2025"
$ws1.Range("C30").Value = "https://codeforces.com/problemset/problem/1777/A"
$ws1.Range("B30").Value = "Everybody likes good arrays"
$ws1.Range("D30").Value = "Easy"
$ws1.Range("E30").Value = "I devised a solution where a group of numbers with the same parity would contribute to the minimum number of operations required.`nThat is, if there is a group of 3 numbers and all are even, then the num ops here would be 2`nSo, we keep on adding this and we get the solution.`nBut this doesnt work !`nThis actually is working but I am facing trouble when the array ends with a group of numbers of the same parity, then I dont go to the else, so the updation is not happening and if I do the updation at last, it creates issues`nAt last, this condition helped:`n    if count > 1:`n        num+=count-1"
$ws1.Range("A30:E30").RowHeight = 187.2

# row 31 (Extremely round, CF1766A)
$ws1.Range("A30:E30").Copy()
$ws1.Range("A31:E31").PasteSpecial(-4122)
$ws1.Range("A31").Value = "30th May, 2025"
$ws1.Range("B31").Value = "Extremely round"
$ws1.Range("C31").Value = "https://codeforces.com/problemset/problem/1766/A"
$ws1.Range("D31").Value = "Easy"
$ws1.Range("E31").Value = "This is a very simple problem, where we need to see the pattern.`nWe observe that the numbers with only one non zero digit appears as :`n1,2,3 .. 9 -> 9 numbers`n10,20,30,….90 -> 9 numbers`n100,200,300,….900 -> 9 numbers`nSo, its observant that after every 10s, there are 9 numbers which are extremely round.`nThus, we keep on dividing our number by 10 until its less than 10, and the number of times we have divided tells us the number of 10s required to build the number and every 10s has 9 extremely round numbers.`nThus the result is number of times divided to make n <10 * 9 + remaining value of n"
$ws1.Range("A31:E31").RowHeight = 187.2

# row 32 (Two permutations, CF1761A)
$ws1.Range("A31:E31").Copy()
$ws1.Range("A32:E32").PasteSpecial(-4122)
$ws1.Range("A32").Value = "1st June, 2025"
$ws1.Range("B32").Value = "Two permutations"
$ws1.Range("C32").Value = "https://codeforces.com/problemset/problem/1761/A"
$ws1.Range("D32").Value = "Easy"
$ws1.Range("E32").Value = "Coming up with the edge cases was a bit difficult. But else the problem is easy.`n<--- a ----> x <----- b ----->`nUpon observation, we can see that if n==a and n==b, i.e, the length of the longest suffix and prefix is equal to n, then the same permutation works for both p and q and thus the ansewr is Yes`nBut, if a+b>n, i.e, a and b together combine to form n, then there is an overlap between a and b and that would imply that would mean :`n<-------a---------->`n                  <------------------------b--->`n<-----------------n---------------------->`nSo, as you can see longest suffix and prefix would essentially be n and hence it would not be possible to form this permutation.`nAnd if there is only one element between a and b or no elements between a annd b, then also its not possible.`nElse it is.`nThis is essentially a very observant visualization problme"
$ws1.Range("A32:E32").RowHeight = 302.4

# Final selection / scroll state to match the saved view
$ws1.Activate()
$ws1.Range("E33").Select()
$excel.ActiveWindow.ScrollRow = 32
$ws3.Range("A19").Select()
